# Verify_CreateInvoiceBySupplier.xlsx - update test data sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the shared-string table rebuilds in
# natural (header-row-then-value-row) order.
$ws.Cells.Clear()

# Row 1 - headers (A1:AF1)
$ws.Range("A1").Value = 'Role'
$ws.Range("B1").Value = 'Location'
$ws.Range("C1").Value = 'Item Description '
$ws.Range("D1").Value = 'UNSPSC Code'
$ws.Range("E1").Value = 'Suggested Supplier(s)'
$ws.Range("F1").Value = 'Category Type'
$ws.Range("G1").Value = 'Category'
$ws.Range("H1").Value = 'Sub Category'
$ws.Range("I1").Value = 'Manufacturer Name'
$ws.Range("J1").Value = 'Manufacturer Part Number'
$ws.Range("K1").Value = 'Quantity'
$ws.Range("L1").Value = 'Unit of Measure'
$ws.Range("M1").Value = 'Price '
$ws.Range("N1").Value = 'ChangeType'
$ws.Range("O1").Value = 'SelectCC '
$ws.Range("P1").Value = 'Role1'
$ws.Range("Q1").Value = 'TaxType'
$ws.Range("R1").Value = 'TaxCode'
$ws.Range("S1").Value = 'ItemName'
$ws.Range("T1").Value = 'ExpectedMsg'
$ws.Range("U1").Value = 'ExpectedStatus'
$ws.Range("V1").Value = 'supplier'
$ws.Range("W1").Value = 'Role2'
$ws.Range("X1").Value = 'Uprice'
$ws.Range("Y1").Value = 'Uquantity'
$ws.Range("Z1").Value = 'LeadTime'
$ws.Range("AA1").Value = 'FreightID'
$ws.Range("AB1").Value = 'Comments'
$ws.Range("AC1").Value = 'CurrentView'
$ws.Range("AD1").Value = 'ActiveIndex'
$ws.Range("AE1").Value = 'ReceivingAgentRole'
$ws.Range("AF1").Value = 'PkgSlipNumber'

# Row 2 - values (A2:AF2)
$ws.Range("A2").Value = 'REQUESTOR'
$ws.Range("B2").Value = 'XEEVA -MJ'
$ws.Range("C2").Value = 'REPOFLOR 100 MG'
$ws.Range("D2").Value = 'UNSPSC001'
$ws.Range("E2").Value = 'Sachin Supplier Magna'
$ws.Range("F2").Value = 'INFORMATION TECHNOLOGY'
$ws.Range("G2").Value = 'HARDWARE'
$ws.Range("H2").Value = 'CELL PHONES'
$ws.Range("I2").Value = 'ARMSTRONG'
$ws.Range("J2").Value = 'MPN001'
$ws.Range("K2").Value = '1;2'
$ws.Range("L2").Value = 'EA-EACH;CU-CUBIC'
$ws.Range("M2").Value = '1;10'
$ws.Range("N2").Value = 'headerlevel '
$ws.Range("O2").Value = 'COMCOMERCIAL '
$ws.Range("P2").Value = 'BUYER'
$ws.Range("Q2").Value = 'Test-Test'
$ws.Range("R2").Value = 'Test_usage-test'
$ws.Range("S2").Value = 'Desktops'
$ws.Range("T2").Value = 'ON-HOLD;Waiting for Approval'
$ws.Range("U2").Value = 'Hold;UnHold'
$ws.Range("V2").Value = 'Supplier'
$ws.Range("W2").Value = 'SUPPLIER'
$ws.Range("X2").Value = "'" + '2'
$ws.Range("Y2").Value = "'" + '10'
$ws.Range("Z2").Value = "'" + '10'
$ws.Range("AA2").Value = 'CIP'
$ws.Range("AB2").Value = 'added comments'
$ws.Range("AC2").Value = 'Myview'
$ws.Range("AD2").Value = "'" + '3'
$ws.Range("AE2").Value = 'REQUESTOR_RECEIVING_AGENT'
$ws.Range("AF2").Value = 'PS123'

# Best-fit-ish column widths for the newly populated columns
$ws.Columns.Item(3).ColumnWidth = 16.5
$ws.Columns.Item(4).ColumnWidth = 12.333333333333332
$ws.Columns.Item(5).ColumnWidth = 20.333333333333336
$ws.Columns.Item(6).ColumnWidth = 26.0
$ws.Columns.Item(7).ColumnWidth = 10.5
$ws.Columns.Item(8).ColumnWidth = 11.666666666666666
$ws.Columns.Item(9).ColumnWidth = 18.166666666666668
$ws.Columns.Item(10).ColumnWidth = 24.333333333333336
$ws.Columns.Item(29).ColumnWidth = 11.5
$ws.Columns.Item(30).ColumnWidth = 10.666666666666666
$ws.Columns.Item(31).ColumnWidth = 28.666666666666668
$ws.Columns.Item(32).ColumnWidth = 14.0

# Restore the active selection to match the authored view
$ws.Range("AE8").Select()
